# "Add score shown system" - appends a new block of localized strings
# (EN_GB / ZH_CN pairs) describing the in-game score/report UI to the
# Localization worksheet, rows 61-70, columns A (EN_GB) and B (ZH_CN).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (EN_GB) values, added first, in this exact order -------
$ws.Range("A61").Value = "REPORT STATUS"
$ws.Range("A62").Value = "REPORT AUTHORITY"
$ws.Range("A63").Value = "SCORE CHANGE"
$ws.Range("A64").Value = "NOT CRIMINAL"
$ws.Range("A65").Value = "CRIMINAL"
$ws.Range("A66").Value = "WRONG"
$ws.Range("A67").Value = "CORRECT"
$ws.Range("A68").Value = "COMBO BONUS"

# --- Column B (ZH_CN) values continue, note row 68's translation is --
# --- entered before rows 61-67's translations, matching source order
$ws.Range("B68").Value = "连对加分"
$ws.Range("B61").Value = "举报状态"
$ws.Range("B62").Value = "举报机构"
$ws.Range("B63").Value = "分数变化"
$ws.Range("B64").Value = "未犯罪"
$ws.Range("B65").Value = "犯罪"
$ws.Range("B66").Value = "错误"
$ws.Range("B67").Value = "正确"

# --- Final pair of rows (current scores / current total) -------------
$ws.Range("A69").Value = "CURRENT SCORES"
$ws.Range("B69").Value = "分数纪录"
$ws.Range("A70").Value = "CURRENT TOTAL"
$ws.Range("B70").Value = "当前总分"

# Apply the same wrap-text style ("s=1") used by the rest of the table
$ws.Range("A61:B70").WrapText = $true

# Update the visible selection to reflect the newly added content
$ws.Range("C71").Select()
